$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.295055866241455
$ws.Range("B1").Value = 2.575385093688965
$ws.Range("C1").Value = 2.478209495544434
$ws.Range("D1").Value = 2.273751258850098
$ws.Range("E1").Value = 0.4238124191761017
